# Update countries & provincias Spain
# Refreshes the COVID-19 country table with newer figures and fixes the
# sort order (by total cases, descending) for a few countries whose case
# counts now cross each other (Barein/Nigeria and Cabo Verde/Surinam/Cuba).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 23:35"

# Row 4: Estados Unidos (updated figures)
$ws.Cells.Item(4, 2).Value = 6455179
$ws.Cells.Item(4, 3).Value = 26039
$ws.Cells.Item(4, 4).Value = 3718433
$ws.Cells.Item(4, 5).Value = 2543612
$ws.Cells.Item(4, 7).Value = 316
$ws.Cells.Item(4, 8).Value = 193134

# Row 6: Brasil (updated figures)
$ws.Cells.Item(6, 2).Value = 4136509
$ws.Cells.Item(6, 3).Value = 13509
$ws.Cells.Item(6, 5).Value = 713167
$ws.Cells.Item(6, 7).Value = 410
$ws.Cells.Item(6, 8).Value = 126640

# Row 24: Alemania (updated figures)
$ws.Cells.Item(24, 2).Value = 251723
$ws.Cells.Item(24, 3).Value = 667
$ws.Cells.Item(24, 5).Value = 16114

# Row 29: Israel (updated figures)
$ws.Cells.Item(29, 2).Value = 130644
$ws.Cells.Item(29, 3).Value = 1708
$ws.Cells.Item(29, 4).Value = 102477
$ws.Cells.Item(29, 5).Value = 27148
$ws.Cells.Item(29, 7).Value = 12
$ws.Cells.Item(29, 8).Value = 1019

# Rows 54-55: Barein overtakes Nigeria in total cases, so they swap places
# Row 54: Nigeria -> Barein (updated figures)
$ws.Cells.Item(54, 1).Value = "Barein"
$ws.Cells.Item(54, 2).Value = 55415
$ws.Cells.Item(54, 3).Value = 644
$ws.Cells.Item(54, 4).Value = 50946
$ws.Cells.Item(54, 5).Value = 4270
$ws.Cells.Item(54, 7).Value = 3
$ws.Cells.Item(54, 8).Value = 199

# Row 55: Barein -> Nigeria (figures unchanged, just moved down a row)
$ws.Cells.Item(55, 1).Value = "Nigeria"
$ws.Cells.Item(55, 2).Value = 54905
$ws.Cells.Item(55, 4).Value = 42922
$ws.Cells.Item(55, 5).Value = 10929
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 1054

# Row 81: Costa de Marfil (updated figures)
$ws.Cells.Item(81, 2).Value = 18588
$ws.Cells.Item(81, 3).Value = 116
$ws.Cells.Item(81, 4).Value = 17472
$ws.Cells.Item(81, 5).Value = 997

# Row 105: Mauritania (updated figures)
$ws.Cells.Item(105, 2).Value = 7142
$ws.Cells.Item(105, 3).Value = 8
$ws.Cells.Item(105, 4).Value = 6669
$ws.Cells.Item(105, 5).Value = 313

# Rows 121-123: Cabo Verde overtakes Surinam and Cuba, so the three rotate
# Row 121: Surinam -> Cabo Verde (updated figures)
$ws.Cells.Item(121, 1).Value = "Cabo Verde"
$ws.Cells.Item(121, 2).Value = 4330
$ws.Cells.Item(121, 3).Value = 55
$ws.Cells.Item(121, 4).Value = 3628
$ws.Cells.Item(121, 5).Value = 660
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 42

# Row 122: Cuba -> Surinam (figures unchanged, just moved down a row)
$ws.Cells.Item(122, 1).Value = "Surinam"
$ws.Cells.Item(122, 2).Value = 4320
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 3426
$ws.Cells.Item(122, 5).Value = 817
$ws.Cells.Item(122, 8).Value = 77

# Row 123: Cabo Verde -> Cuba (figures unchanged, just moved down a row)
$ws.Cells.Item(123, 1).Value = "Cuba"
$ws.Cells.Item(123, 2).Value = 4309
$ws.Cells.Item(123, 3).Value = 11
$ws.Cells.Item(123, 4).Value = 3590
$ws.Cells.Item(123, 5).Value = 618
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 101
